$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename reviewer "Tibor" -> "Thea" for the Framingham rows (double-checked) ---
$ws.Range("J4").Value = "Thea"
$ws.Range("J5").Value = "Thea"
$ws.Range("J6").Value = "Thea"

# --- Expand the ARIC model description (A5) into rich text with extra detail ---
$ws.Range("A5").Value = "ARIC Model - Clinical variables plus fasting glucose and lipids (Schmidt et al.)"

# "ARIC" (chars 1-4) keeps the default cell font (no explicit run formatting)
# " Model - Clinical variables plus fasting glucose and lipids " (chars 5-64)
$run2 = $ws.Range("A5").Characters(5, 60)
$run2.Font.Size = 9
$run2.Font.Name = "Calibri"

# "(Schmidt et al.)" (chars 65-80)
$run3 = $ws.Range("A5").Characters(65, 16)
$run3.Font.Size = 9
$run3.Font.Name = "Calibri"

# --- Update the sheet view: scroll + selection moved to around J9 ---
$ws.Activate()
$ws.Range("J9").Select()
